$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.194401264190674
$ws.Range("B1").Value = 2.128674268722534
$ws.Range("C1").Value = 3.880314350128174
$ws.Range("D1").Value = 3.30093240737915
$ws.Range("E1").Value = 1.13032066822052
